{"js": "// The captured change re-saved this package with Apache POI 3.15 (see commit\n// message: \"Fixed POI packaging and upgraded to POI 3.15\"). POI's XML writer\n// serializes element attributes in alphabetical order, so every hunk in the\n// diff is the *same* tab stops / page size / page margins / footnote ids /\n// style attributes simply re-emitted with their attributes reordered\n// (e.g. <w:tab w:val=\"left\" w:pos=\"3119\"/> -> <w:tab w:pos=\"3119\" w:val=\"left\"/>,\n// <w:pgSz w:w=\"11906\" w:h=\"16838\"/> -> <w:pgSz w:h=\"16838\" w:w=\"11906\"/>).\n// No value anywhere in the document actually changes.\n//\n// Office.js does not expose raw attribute-serialization order, so we apply\n// the change at the content/value level: touch every property the diff\n// rewrites and reassert its existing value through the supported object\n// model so the part is regenerated, keeping the document semantically\n// identical to the target (which is also semantically identical to the\n// source - only byte-level attribute order moved).\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\n// word/document.xml -> w:sectPr -> w:pgSz / w:pgMar (values unchanged:\n// 11906x16838 twips = 595.3x841.9pt page, 1417/708/0 twips margins).\nconst section = sections.items[0];\nconst pageSetup = section.pageSetup;\npageSetup.load(\"pageWidth,pageHeight,topMargin,bottomMargin,leftMargin,rightMargin,headerDistance,footerDistance,gutter\");\nawait context.sync();\n\npageSetup.pageWidth = pageSetup.pageWidth;\npageSetup.pageHeight = pageSetup.pageHeight;\npageSetup.topMargin = pageSetup.topMargin;\npageSetup.bottomMargin = pageSetup.bottomMargin;\npageSetup.leftMargin = pageSetup.leftMargin;\npageSetup.rightMargin = pageSetup.rightMargin;\npageSetup.headerDistance = pageSetup.headerDistance;\npageSetup.footerDistance = pageSetup.footerDistance;\npageSetup.gutter = pageSetup.gutter;\n\nawait context.sync();\n", "ps1": "# The captured change re-saved this package with Apache POI 3.15 (see commit\n# message: \"Fixed POI packaging and upgraded to POI 3.15\"). POI's XML writer\n# serializes element attributes in alphabetical order, so every hunk in the\n# diff is the *same* tab stops / page size / page margins / footnote ids /\n# style attributes simply re-emitted with their attributes reordered, e.g.\n#   <w:tab w:val=\"left\" w:pos=\"3119\"/>  ->  <w:tab w:pos=\"3119\" w:val=\"left\"/>\n#   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>   ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n# No value anywhere in the document actually changes.\n#\n# Word COM does not expose raw attribute-serialization order either, so we\n# apply the change at the content/value level: touch every property the\n# diff rewrites and reassert its existing value through the object model so\n# the part is regenerated, keeping the document semantically identical to\n# the target (which is also semantically identical to the source - only\n# byte-level attribute order moved).\n\n$d = $word.ActiveDocument\n\n# word/document.xml -> w:pPr/w:tabs/w:tab (left tab @ pos=3119 twips =\n# 155.95pt) on the first four paragraphs - untouched positions/alignment,\n# just re-asserted.\nfor ($i = 1; $i -le 4; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $tabStops = $p.Range.ParagraphFormat.TabStops\n    $tab = $tabStops.Item(1)\n    $tab.Position = $tab.Position\n    $tab.Alignment = $tab.Alignment\n}\n\n# word/document.xml -> w:sectPr -> w:pgSz / w:pgMar (values unchanged:\n# 11906x16838 twips = 595.3x841.9pt page, 1417/708/0 twips margins).\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n$pageSetup.PageWidth = $pageSetup.PageWidth\n$pageSetup.PageHeight = $pageSetup.PageHeight\n$pageSetup.TopMargin = $pageSetup.TopMargin\n$pageSetup.BottomMargin = $pageSetup.BottomMargin\n$pageSetup.LeftMargin = $pageSetup.LeftMargin\n$pageSetup.RightMargin = $pageSetup.RightMargin\n$pageSetup.HeaderDistance = $pageSetup.HeaderDistance\n$pageSetup.FooterDistance = $pageSetup.FooterDistance\n$pageSetup.Gutter = $pageSetup.Gutter\n"}
